$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: RandomForestRegressor - name unchanged, values updated
$ws.Range("B3").Value = 0.9781751630831635
$ws.Range("C3").Value = 0.9782633397664237
$ws.Range("D3").Value = 0.9774636623307407

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9788673185656463
$ws.Range("C4").Value = 0.9788254521393971
$ws.Range("D4").Value = 0.9786215090977686

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8220022803819297
$ws.Range("C5").Value = 0.8388643957271017
$ws.Range("D5").Value = 0.8139862639266873
